$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 88.75
$ws.Range("H29").Value = 1742.619
$ws.Range("J29").Value = 2236.875
$ws.Range("L29").Value = 6710.625
$ws.Range("N29").Value = -7272.625
$ws.Range("H31").Value = 980000
$ws.Range("I31").Value = 980000
$ws.Range("K31").Value = 2940000
$ws.Range("M31").Value = -2939770
$ws.Range("H32").Value = 1248
$ws.Range("H33").Value = 265.89285
$ws.Range("I33").Value = 214.09091
$ws.Range("J33").Value = 455.83334
$ws.Range("K33").Value = 214.09091
$ws.Range("L33").Value = 455.83334
$ws.Range("M33").Value = 14.90908999999999
$ws.Range("N33").Value = -913.83334
$ws.Range("H38").Value = 1618.2195
$ws.Range("I38").Value = 217.875
$ws.Range("J38").Value = 1957.697
$ws.Range("K38").Value = 653.625
$ws.Range("L38").Value = 5873.090999999999
$ws.Range("M38").Value = -281.625
$ws.Range("N38").Value = -6617.090999999999
$ws.Range("H41").Value = 1837.1052
$ws.Range("I41").Value = 1750.3846
$ws.Range("J41").Value = 2025
$ws.Range("K41").Value = 1750.3846
$ws.Range("L41").Value = 2025
$ws.Range("M41").Value = -1310.3846
$ws.Range("N41").Value = -2905
$ws.Range("H47").Value = 3000
$ws.Range("I47").Value = 3000
$ws.Range("K47").Value = 3000
$ws.Range("M47").Value = -2028
$ws.Range("H112").Value = 3483.5264
$ws.Range("I112").Value = 1025
$ws.Range("J112").Value = 4139.1333
$ws.Range("K112").Value = 3075
$ws.Range("L112").Value = 12417.3999
$ws.Range("M112").Value = -1967
$ws.Range("N112").Value = -14633.3999
$ws.Range("H132").Value = 8554017
$ws.Range("I132").Value = 13340238
$ws.Range("J132").Value = 7191.857
$ws.Range("K132").Value = 40020714
$ws.Range("L132").Value = 21575.571
$ws.Range("M132").Value = -40018184
$ws.Range("N132").Value = -26635.571
$ws.Range("H137").Value = 1495.3556
$ws.Range("I137").Value = 1252.125
$ws.Range("J137").Value = 1773.3334
$ws.Range("K137").Value = 3756.375
$ws.Range("L137").Value = 5320.0002
$ws.Range("M137").Value = -1206.375
$ws.Range("N137").Value = -10420.0002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5693.2705
$ws.Range("I32").Value = 6077.9707
$ws.Range("J32").Value = 1333.3334
$ws.Range("K32").Value = 6077.9707
$ws.Range("L32").Value = 1333.3334
$ws.Range("M32").Value = -5790.9707
$ws.Range("N32").Value = -1907.3334
$ws.Range("H61").Value = 90910340
$ws.Range("I61").Value = 142857950
$ws.Range("J61").Value = 2004.5
$ws.Range("K61").Value = 142857950
$ws.Range("L61").Value = 2004.5
$ws.Range("M61").Value = -142857738
$ws.Range("N61").Value = -2428.5
$ws.Range("H63").Value = 2050.4888
$ws.Range("I63").Value = 1910.7333
$ws.Range("K63").Value = 1910.7333
$ws.Range("M63").Value = -1224.7333
$ws.Range("H66").Value = 2050.4888
$ws.Range("I66").Value = 1910.7333
$ws.Range("K66").Value = 9553.666500000001
$ws.Range("M66").Value = -6121.666500000001
$ws.Range("H74").Value = 2552.3333
$ws.Range("I74").Value = 2033.4286
$ws.Range("J74").Value = 3278.8
$ws.Range("K74").Value = 2033.4286
$ws.Range("L74").Value = 3278.8
$ws.Range("M74").Value = -1159.4286
$ws.Range("N74").Value = -5026.8
$ws.Range("H77").Value = 2552.3333
$ws.Range("I77").Value = 2033.4286
$ws.Range("J77").Value = 3278.8
$ws.Range("K77").Value = 10167.143
$ws.Range("L77").Value = 16394
$ws.Range("M77").Value = -5799.143
$ws.Range("N77").Value = -25130
$ws.Range("H102").Value = 11907463
$ws.Range("I102").Value = 18521054
$ws.Range("K102").Value = 18521054
$ws.Range("M102").Value = -18519432
$ws.Range("H132").Value = 3288.5
$ws.Range("I132").Value = 2842
$ws.Range("K132").Value = 8526
$ws.Range("M132").Value = -5996
$ws.Range("H136").Value = 90910340
$ws.Range("I136").Value = 142857950
$ws.Range("J136").Value = 2004.5
$ws.Range("K136").Value = 428573850
$ws.Range("L136").Value = 6013.5
$ws.Range("M136").Value = -428571300
$ws.Range("N136").Value = -11113.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 63119290
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
$ws.Range("H107").Value = 1364.5
$ws.Range("I107").Value = 986.1
$ws.Range("J107").Value = 3256.5
$ws.Range("K107").Value = 986.1
$ws.Range("L107").Value = 3256.5
$ws.Range("M107").Value = 933.9
$ws.Range("N107").Value = -7096.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 100131.57
$ws.Range("I22").Value = 144.2
$ws.Range("J22").Value = 350100
$ws.Range("K22").Value = 144.2
$ws.Range("L22").Value = 350100
$ws.Range("M22").Value = 205.8
$ws.Range("N22").Value = -350800
$ws.Range("H31").Value = 1515.5
$ws.Range("I31").Value = 1497.0769
$ws.Range("J31").Value = 1549.7142
$ws.Range("K31").Value = 1497.0769
$ws.Range("L31").Value = 1549.7142
$ws.Range("M31").Value = -1202.0769
$ws.Range("N31").Value = -2139.7142
$ws.Range("H34").Value = 1515.5
$ws.Range("I34").Value = 1497.0769
$ws.Range("J34").Value = 1549.7142
$ws.Range("K34").Value = 1497.0769
$ws.Range("L34").Value = 1549.7142
$ws.Range("M34").Value = -1295.0769
$ws.Range("N34").Value = -1953.7142
$ws.Range("H99").Value = 1330.4736
$ws.Range("I99").Value = 1256.6923
$ws.Range("J99").Value = 1490.3334
$ws.Range("K99").Value = 1256.6923
$ws.Range("L99").Value = 1490.3334
$ws.Range("M99").Value = 241.3077000000001
$ws.Range("N99").Value = -4486.3334
$ws.Range("H126").Value = 1330.4736
$ws.Range("I126").Value = 1256.6923
$ws.Range("J126").Value = 1490.3334
$ws.Range("K126").Value = 3770.0769
$ws.Range("L126").Value = 4471.0002
$ws.Range("M126").Value = -1300.0769
$ws.Range("N126").Value = -9411.0002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 661.931
$ws.Range("I113").Value = 578.25
$ws.Range("J113").Value = 693.8095
$ws.Range("K113").Value = 1734.75
$ws.Range("L113").Value = 2081.4285
$ws.Range("M113").Value = 435.25
$ws.Range("N113").Value = -6421.4285
$ws.Range("H131").Value = 19611416
$ws.Range("I131").Value = 142857460
$ws.Range("J131").Value = 4092.3865
$ws.Range("K131").Value = 428572380
$ws.Range("L131").Value = 12277.1595
$ws.Range("M131").Value = -428567340
$ws.Range("N131").Value = -22357.1595
$ws.Range("H138").Value = 1470.1714
$ws.Range("I138").Value = 912.0909
$ws.Range("J138").Value = 2414.6155
$ws.Range("K138").Value = 2736.2727
$ws.Range("L138").Value = 7243.8465
$ws.Range("M138").Value = 2403.7273
$ws.Range("N138").Value = -17523.8465
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4253.5
$ws.Range("I80").Value = 3200.8333
$ws.Range("J80").Value = 5832.5
$ws.Range("K80").Value = 3200.8333
$ws.Range("L80").Value = 5832.5
$ws.Range("M80").Value = -2202.8333
$ws.Range("N80").Value = -7828.5
$ws.Range("H83").Value = 4253.5
$ws.Range("I83").Value = 3200.8333
$ws.Range("J83").Value = 5832.5
$ws.Range("K83").Value = 16004.1665
$ws.Range("L83").Value = 29162.5
$ws.Range("M83").Value = -11012.1665
$ws.Range("N83").Value = -39146.5
$ws.Range("H132").Value = 2974.3809
$ws.Range("I132").Value = 2673.5293
$ws.Range("K132").Value = 8020.5879
$ws.Range("M132").Value = -5490.5879
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3764.6428
$ws.Range("I46").Value = 554.5
$ws.Range("J46").Value = 6172.25
$ws.Range("K46").Value = 554.5
$ws.Range("L46").Value = 6172.25
$ws.Range("M46").Value = -366.5
$ws.Range("N46").Value = -6548.25
$ws.Range("H124").Value = 39500
$ws.Range("J124").Value = 39500
$ws.Range("L124").Value = 39500
$ws.Range("N124").Value = -49320
